# Remove the stray space that follows the literal "<" right before the
# "https://github.com/tbrowder" link, in every copy of the handout found
# in this document (the link is repeated once per handout block).
#
# Before: "< https://github.com/tbrowder>"
# After:  "<https://github.com/tbrowder>"

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $rngText = $p.Range.Text
    $pos = $rngText.IndexOf("< https://github.com")
    if ($pos -ge 0) {
        $paraStart = $p.Range.Start
        # The offending character is the single space immediately after '<'.
        $spaceStart = $paraStart + $pos + 1
        $spaceEnd = $spaceStart + 1
        $spaceRange = $d.Range($spaceStart, $spaceEnd)
        $spaceRange.Delete()
    }
}
